# Auto-generated Excel COM-interop script to apply the weekly Coco price-sheet update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows 3-27: refresh Fecha/Volumen/Precios/Precio-Kg columns ---
$ws.Cells.Item(3, 4).Value = 44425
$ws.Cells.Item(3, 13).Value = 15
$ws.Cells.Item(3, 14).Value = 24000
$ws.Cells.Item(3, 15).Value = 24000
$ws.Cells.Item(3, 16).Value = 24000
$ws.Cells.Item(3, 19).Value = 1200

$ws.Cells.Item(4, 4).Value = 44392
$ws.Cells.Item(4, 13).Value = 10
$ws.Cells.Item(4, 14).Value = 24000
$ws.Cells.Item(4, 15).Value = 24000
$ws.Cells.Item(4, 16).Value = 24000
$ws.Cells.Item(4, 19).Value = 1200

$ws.Cells.Item(5, 4).Value = 44249
$ws.Cells.Item(5, 13).Value = 15
$ws.Cells.Item(5, 14).Value = 25000
$ws.Cells.Item(5, 15).Value = 25000
$ws.Cells.Item(5, 16).Value = 25000
$ws.Cells.Item(5, 19).Value = 1250

$ws.Cells.Item(6, 4).Value = 44424
$ws.Cells.Item(6, 13).Value = 25
$ws.Cells.Item(6, 14).Value = 24000
$ws.Cells.Item(6, 15).Value = 24000
$ws.Cells.Item(6, 16).Value = 24000
$ws.Cells.Item(6, 19).Value = 1200

$ws.Cells.Item(7, 4).Value = 44414
$ws.Cells.Item(7, 13).Value = 15
$ws.Cells.Item(7, 14).Value = 25000
$ws.Cells.Item(7, 15).Value = 25000
$ws.Cells.Item(7, 16).Value = 25000
$ws.Cells.Item(7, 19).Value = 1250

$ws.Cells.Item(8, 4).Value = 44356
$ws.Cells.Item(8, 13).Value = 15
$ws.Cells.Item(8, 14).Value = 24000
$ws.Cells.Item(8, 15).Value = 24000
$ws.Cells.Item(8, 16).Value = 24000
$ws.Cells.Item(8, 19).Value = 1200

$ws.Cells.Item(9, 4).Value = 44396
$ws.Cells.Item(9, 13).Value = 12
$ws.Cells.Item(9, 14).Value = 24000
$ws.Cells.Item(9, 15).Value = 24000
$ws.Cells.Item(9, 16).Value = 24000
$ws.Cells.Item(9, 19).Value = 1200

$ws.Cells.Item(10, 4).Value = 44221
$ws.Cells.Item(10, 13).Value = 30
$ws.Cells.Item(10, 14).Value = 25000
$ws.Cells.Item(10, 15).Value = 25000
$ws.Cells.Item(10, 16).Value = 25000
$ws.Cells.Item(10, 19).Value = 1250

$ws.Cells.Item(11, 4).Value = 44175
$ws.Cells.Item(11, 13).Value = 25
$ws.Cells.Item(11, 14).Value = 23000
$ws.Cells.Item(11, 15).Value = 23000
$ws.Cells.Item(11, 16).Value = 23000
$ws.Cells.Item(11, 19).Value = 1150

$ws.Cells.Item(12, 4).Value = 44363
$ws.Cells.Item(12, 13).Value = 30
$ws.Cells.Item(12, 14).Value = 24000
$ws.Cells.Item(12, 15).Value = 24000
$ws.Cells.Item(12, 16).Value = 24000
$ws.Cells.Item(12, 19).Value = 1200

$ws.Cells.Item(13, 4).Value = 44349
$ws.Cells.Item(13, 13).Value = 30
$ws.Cells.Item(13, 14).Value = 24000
$ws.Cells.Item(13, 15).Value = 24000
$ws.Cells.Item(13, 16).Value = 24000
$ws.Cells.Item(13, 19).Value = 1200

$ws.Cells.Item(14, 4).Value = 44421
$ws.Cells.Item(14, 13).Value = 20
$ws.Cells.Item(14, 14).Value = 24000
$ws.Cells.Item(14, 15).Value = 24000
$ws.Cells.Item(14, 16).Value = 24000
$ws.Cells.Item(14, 19).Value = 1200

$ws.Cells.Item(15, 4).Value = 44222
$ws.Cells.Item(15, 13).Value = 15
$ws.Cells.Item(15, 14).Value = 25000
$ws.Cells.Item(15, 15).Value = 25000
$ws.Cells.Item(15, 16).Value = 25000
$ws.Cells.Item(15, 19).Value = 1250

$ws.Cells.Item(16, 4).Value = 44377
$ws.Cells.Item(16, 13).Value = 15
$ws.Cells.Item(16, 14).Value = 20000
$ws.Cells.Item(16, 15).Value = 20000
$ws.Cells.Item(16, 16).Value = 20000
$ws.Cells.Item(16, 19).Value = 1000

$ws.Cells.Item(17, 4).Value = 44400
$ws.Cells.Item(17, 13).Value = 5
$ws.Cells.Item(17, 14).Value = 24000
$ws.Cells.Item(17, 15).Value = 24000
$ws.Cells.Item(17, 16).Value = 24000
$ws.Cells.Item(17, 19).Value = 1200

$ws.Cells.Item(18, 4).Value = 44426
$ws.Cells.Item(18, 13).Value = 15
$ws.Cells.Item(18, 14).Value = 24000
$ws.Cells.Item(18, 15).Value = 24000
$ws.Cells.Item(18, 16).Value = 24000
$ws.Cells.Item(18, 19).Value = 1200

$ws.Cells.Item(19, 4).Value = 44194
$ws.Cells.Item(19, 13).Value = 20
$ws.Cells.Item(19, 14).Value = 20000
$ws.Cells.Item(19, 15).Value = 20000
$ws.Cells.Item(19, 16).Value = 20000
$ws.Cells.Item(19, 19).Value = 1000

$ws.Cells.Item(20, 4).Value = 44390
$ws.Cells.Item(20, 13).Value = 10
$ws.Cells.Item(20, 14).Value = 24000
$ws.Cells.Item(20, 15).Value = 24000
$ws.Cells.Item(20, 16).Value = 24000
$ws.Cells.Item(20, 19).Value = 1200

$ws.Cells.Item(21, 4).Value = 44412
$ws.Cells.Item(21, 13).Value = 20
$ws.Cells.Item(21, 14).Value = 25000
$ws.Cells.Item(21, 15).Value = 25000
$ws.Cells.Item(21, 16).Value = 25000
$ws.Cells.Item(21, 19).Value = 1250

$ws.Cells.Item(22, 4).Value = 44214
$ws.Cells.Item(22, 13).Value = 15
$ws.Cells.Item(22, 14).Value = 25000
$ws.Cells.Item(22, 15).Value = 25000
$ws.Cells.Item(22, 16).Value = 25000
$ws.Cells.Item(22, 19).Value = 1250

$ws.Cells.Item(23, 4).Value = 44238
$ws.Cells.Item(23, 13).Value = 30
$ws.Cells.Item(23, 14).Value = 25000
$ws.Cells.Item(23, 15).Value = 25000
$ws.Cells.Item(23, 16).Value = 25000
$ws.Cells.Item(23, 19).Value = 1250

$ws.Cells.Item(24, 4).Value = 44231
$ws.Cells.Item(24, 13).Value = 15
$ws.Cells.Item(24, 14).Value = 25000
$ws.Cells.Item(24, 15).Value = 25000
$ws.Cells.Item(24, 16).Value = 25000
$ws.Cells.Item(24, 19).Value = 1250

$ws.Cells.Item(25, 4).Value = 44391
$ws.Cells.Item(25, 13).Value = 10
$ws.Cells.Item(25, 14).Value = 24000
$ws.Cells.Item(25, 15).Value = 24000
$ws.Cells.Item(25, 16).Value = 24000
$ws.Cells.Item(25, 19).Value = 1200

$ws.Cells.Item(26, 4).Value = 44389
$ws.Cells.Item(26, 13).Value = 20
$ws.Cells.Item(26, 14).Value = 24000
$ws.Cells.Item(26, 15).Value = 24000
$ws.Cells.Item(26, 16).Value = 24000
$ws.Cells.Item(26, 19).Value = 1200

$ws.Cells.Item(27, 4).Value = 44251
$ws.Cells.Item(27, 13).Value = 15
$ws.Cells.Item(27, 14).Value = 25000
$ws.Cells.Item(27, 15).Value = 25000
$ws.Cells.Item(27, 16).Value = 25000
$ws.Cells.Item(27, 19).Value = 1250

# --- Append new weekly data rows 28-34 ---
$newRows = @(
    @{ D = 44419; M = 40; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ D = 44420; M = 35; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ D = 44382; M = 15; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ D = 44232; M = 15; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ D = 44398; M = 15; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ D = 44334; M = 20; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ D = 44418; M = 20; N = 24000; O = 24000; P = 24000; S = 1200 }
)

$r = 28
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108007
    $ws.Cells.Item($r, 10).Value = "Coco"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = "Primera"
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = "$/malla 20 unidades"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = 20
    $r++
}

